# Scheduled-runner price refresh for the Phantom_Profits workbook.
# Re-pulls Universalis market data (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ],
# LeveProfit[NQ/HQ]) for a set of leves across all Disciple of the Hand sheets and
# writes the refreshed cached values back into columns H:N.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 416.2143
$ws.Range("I33").Value = 403.54544
$ws.Range("K33").Value = 403.54544
$ws.Range("M33").Value = -174.54544
$ws.Range("H58").Value = 763.75
$ws.Range("I58").Value = 527.5
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1582.5
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1432.5
$ws.Range("N58").Value = -3300
$ws.Range("H96").Value = 1667.6154
$ws.Range("J96").Value = 1261.7142
$ws.Range("L96").Value = 3785.1426
$ws.Range("N96").Value = -6531.142599999999
$ws.Range("H98").Value = 935.1905
$ws.Range("I98").Value = 297.35294
$ws.Range("K98").Value = 297.35294
$ws.Range("M98").Value = 1200.64706
$ws.Range("H100").Value = 1927.7
$ws.Range("I100").Value = 1572.25
$ws.Range("K100").Value = 1572.25
$ws.Range("M100").Value = -1031.25
$ws.Range("H122").Value = 935.1905
$ws.Range("I122").Value = 297.35294
$ws.Range("K122").Value = 892.05882
$ws.Range("M122").Value = 1557.94118
$ws.Range("H132").Value = 3949.9473
$ws.Range("I132").Value = 3908.7646
$ws.Range("J132").Value = 4300
$ws.Range("K132").Value = 11726.2938
$ws.Range("L132").Value = 12900
$ws.Range("M132").Value = -9196.293799999999
$ws.Range("N132").Value = -17960
$ws.Range("H137").Value = 1996.1428
$ws.Range("J137").Value = 1621.25
$ws.Range("L137").Value = 4863.75
$ws.Range("N137").Value = -9963.75
$ws.Range("H138").Value = 3298
$ws.Range("I138").Value = 700
$ws.Range("K138").Value = 2100
$ws.Range("M138").Value = 3040

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1098.3334
$ws.Range("I2").Value = 1098.3334
$ws.Range("K2").Value = 1098.3334
$ws.Range("M2").Value = -985.3334
$ws.Range("H5").Value = 1289.3
$ws.Range("I5").Value = 2769.75
$ws.Range("K5").Value = 2769.75
$ws.Range("M5").Value = -2657.75
$ws.Range("H32").Value = 1027.3226
$ws.Range("I32").Value = 925.7586
$ws.Range("K32").Value = 925.7586
$ws.Range("M32").Value = -638.7586
$ws.Range("H63").Value = 6910.684
$ws.Range("I63").Value = 3715.6155
$ws.Range("J63").Value = 13833.333
$ws.Range("K63").Value = 3715.6155
$ws.Range("L63").Value = 13833.333
$ws.Range("M63").Value = -3029.6155
$ws.Range("N63").Value = -15205.333
$ws.Range("H66").Value = 6910.684
$ws.Range("I66").Value = 3715.6155
$ws.Range("J66").Value = 13833.333
$ws.Range("K66").Value = 18578.0775
$ws.Range("L66").Value = 69166.66500000001
$ws.Range("M66").Value = -15146.0775
$ws.Range("N66").Value = -76030.66500000001
$ws.Range("H74").Value = 1181.375
$ws.Range("I74").Value = 1060.2
$ws.Range("K74").Value = 1060.2
$ws.Range("M74").Value = -186.2
$ws.Range("H77").Value = 1181.375
$ws.Range("I77").Value = 1060.2
$ws.Range("K77").Value = 5301
$ws.Range("M77").Value = -933
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H116").Value = 1098.3334
$ws.Range("I116").Value = 1098.3334
$ws.Range("K116").Value = 1098.3334
$ws.Range("M116").Value = 1195.6666
$ws.Range("H132").Value = 3181.05
$ws.Range("J132").Value = 4721.8887
$ws.Range("L132").Value = 14165.6661
$ws.Range("N132").Value = -19225.6661

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1098.3334
$ws.Range("I3").Value = 1098.3334
$ws.Range("K3").Value = 1098.3334
$ws.Range("M3").Value = -984.3334
$ws.Range("H4").Value = 1289.3
$ws.Range("I4").Value = 2769.75
$ws.Range("K4").Value = 2769.75
$ws.Range("M4").Value = -2654.75
$ws.Range("H99").Value = 5496669.5
$ws.Range("I99").Value = 6411581
$ws.Range("K99").Value = 6411581
$ws.Range("M99").Value = -6410083
$ws.Range("H105").Value = 100003970
$ws.Range("I105").Value = 2135
$ws.Range("K105").Value = 2135
$ws.Range("M105").Value = -388
$ws.Range("H107").Value = 1830
$ws.Range("I107").Value = 1830
$ws.Range("K107").Value = 1830
$ws.Range("M107").Value = 90

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.15385000000001
$ws.Range("I7").Value = 98.5
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 98.5
$ws.Range("L7").Value = 16
$ws.Range("M7").Value = 14.5
$ws.Range("N7").Value = -242
$ws.Range("H22").Value = 5863.56
$ws.Range("I22").Value = 1289.4166
$ws.Range("J22").Value = 10085.846
$ws.Range("K22").Value = 1289.4166
$ws.Range("L22").Value = 10085.846
$ws.Range("M22").Value = -939.4166
$ws.Range("N22").Value = -10785.846
$ws.Range("H63").Value = 85136
$ws.Range("J63").Value = 85136
$ws.Range("L63").Value = 85136
$ws.Range("N63").Value = -86508
$ws.Range("H64").Value = 50271
$ws.Range("J64").Value = 50271
$ws.Range("L64").Value = 50271
$ws.Range("N64").Value = -50767
$ws.Range("H66").Value = 85136
$ws.Range("J66").Value = 85136
$ws.Range("L66").Value = 255408
$ws.Range("N66").Value = -262272
$ws.Range("H67").Value = 50271
$ws.Range("J67").Value = 50271
$ws.Range("L67").Value = 50271
$ws.Range("N67").Value = -51987
$ws.Range("H88").Value = 13749
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 13749
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 30
$ws.Range("K3").Value = 90
$ws.Range("M3").Value = 22
$ws.Range("H108").Value = 1486
$ws.Range("I108").Value = 1486
$ws.Range("K108").Value = 4458
$ws.Range("M108").Value = -1578

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2151.9333
$ws.Range("I122").Value = 1652.6364
$ws.Range("J122").Value = 3525
$ws.Range("K122").Value = 4957.9092
$ws.Range("L122").Value = 10575
$ws.Range("M122").Value = -2507.9092
$ws.Range("N122").Value = -15475
$ws.Range("H132").Value = 5499.8335
$ws.Range("I132").Value = 4999.75
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 14999.25
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -12469.25
$ws.Range("N132").Value = -24560

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H22").Value = 2125.2307
$ws.Range("J22").Value = 2717.6
$ws.Range("L22").Value = 2717.6
$ws.Range("N22").Value = -3307.6
$ws.Range("H27").Value = 2125.2307
$ws.Range("J27").Value = 2717.6
$ws.Range("L27").Value = 2717.6
$ws.Range("N27").Value = -2931.6
$ws.Range("H61").Value = 1533.0869
$ws.Range("I61").Value = 1362.1666
$ws.Range("K61").Value = 1362.1666
$ws.Range("M61").Value = -1160.1666
$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 10000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -8877
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 30000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -24384
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 2999.5
$ws.Range("I100").Value = 2999.5
$ws.Range("K100").Value = 2999.5
$ws.Range("M100").Value = -2458.5
$ws.Range("H113").Value = 1533.0869
$ws.Range("I113").Value = 1362.1666
$ws.Range("K113").Value = 1362.1666
$ws.Range("M113").Value = 807.8334
$ws.Range("H122").Value = 3966
$ws.Range("I122").Value = 3949.5
$ws.Range("K122").Value = 11848.5
$ws.Range("M122").Value = -9398.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 2869.1052
$ws.Range("I132").Value = 2770.6428
$ws.Range("J132").Value = 3144.8
$ws.Range("K132").Value = 8311.928400000001
$ws.Range("L132").Value = 9434.400000000001
$ws.Range("M132").Value = -5781.928400000001
$ws.Range("N132").Value = -14494.4
$ws.Range("H136").Value = 8383.143
$ws.Range("I136").Value = 10937.4
$ws.Range("J136").Value = 1997.5
$ws.Range("K136").Value = 32812.2
$ws.Range("L136").Value = 5992.5
$ws.Range("M136").Value = -30262.2
$ws.Range("N136").Value = -11092.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 27904.285
$ws.Range("J69").Value = 27904.285
$ws.Range("L69").Value = 27904.285
$ws.Range("N69").Value = -29402.285
$ws.Range("H72").Value = 27904.285
$ws.Range("J72").Value = 27904.285
$ws.Range("L72").Value = 83712.855
$ws.Range("N72").Value = -91200.855
$ws.Range("H122").Value = 1497.5
$ws.Range("I122").Value = 1497.5
$ws.Range("K122").Value = 4492.5
$ws.Range("M122").Value = -2042.5
$ws.Range("H130").Value = 24500
$ws.Range("J130").Value = 24500
$ws.Range("L130").Value = 24500
$ws.Range("N130").Value = -34540

Write-Output "Applied 244 cell updates across 8 sheets"
